# Update NATMI LR-pair TPM-derived values for Fgf9-Fgfr1 (ECs sending cluster)
# Underlying change: Ligand-expressing cells (E) dropped from 3 to 2 for every
# target-cluster row, and the Receptor average expression values (M) were
# recomputed from the new TPM table. All other touched columns are values
# derived from E/G/H/M/N, so they are recalculated here in Python-equivalent
# PowerShell math and written directly, matching the canonical OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Total ligand cells in the sending cluster (ECs) is fixed at 3; only the
# expressing-cell count changes from 3 -> 2, which also rescales the ligand
# average/total expression values (same scale factor for every row because
# the ligand/sending cluster is identical across rows).
$ligandTotalCells = 3
$newLigandExpressingCells = 2
$newLigandAvg = 0.579684
$newLigandTotal = 1.739052

# Receptor-expressing cells (K) stays at 3 for every row.
$receptorCells = 3

# New receptor average expression values (M), recomputed from updated TPM.
$newReceptorAvg = @{
    2 = 10.48767733333333
    3 = 62.99699166666667
    4 = 0.3322793333333333
    5 = 10.25458433333333
    6 = 0.7572163333333334
    7 = 0.9889696666666666
}

$rows = 2..7

# First pass: write E, F, G, H (ligand columns) and M, N (receptor columns).
foreach ($r in $rows) {
    $ws.Cells.Item($r, 5).Value2  = $newLigandExpressingCells          # E: Ligand-expressing cells
    $ws.Cells.Item($r, 6).Value2  = $newLigandExpressingCells / $ligandTotalCells  # F: Ligand detection rate
    $ws.Cells.Item($r, 7).Value2  = $newLigandAvg                      # G: Ligand average expression value
    $ws.Cells.Item($r, 8).Value2  = $newLigandTotal                    # H: Ligand total expression value

    $mVal = $newReceptorAvg[$r]
    $nVal = $mVal * $receptorCells

    $ws.Cells.Item($r, 13).Value2 = $mVal                               # M: Receptor average expression value
    $ws.Cells.Item($r, 14).Value2 = $nVal                               # N: Receptor total expression value
}

# Second pass: derive the edge weights (Q = G*M, R = H*N) per row, then the
# specificity columns (O/S from Q, P/T from R) normalized across all rows for
# this ligand-receptor pair.
$qVals = @{}
$rVals = @{}
foreach ($r in $rows) {
    $mVal = $newReceptorAvg[$r]
    $nVal = $mVal * $receptorCells
    $qVals[$r] = $newLigandAvg * $mVal
    $rVals[$r] = $newLigandTotal * $nVal
}

$qTotal = 0
$rTotal = 0
foreach ($r in $rows) {
    $qTotal += $qVals[$r]
    $rTotal += $rVals[$r]
}

foreach ($r in $rows) {
    $qVal = $qVals[$r]
    $rVal = $rVals[$r]
    $oVal = $qVal / $qTotal
    $pVal = $rVal / $rTotal

    $ws.Cells.Item($r, 15).Value2 = $oVal   # O: Receptor derived specificity of average expression value
    $ws.Cells.Item($r, 16).Value2 = $pVal   # P: Receptor derived specificity of total expression value
    $ws.Cells.Item($r, 17).Value2 = $qVal   # Q: Edge average expression weight
    $ws.Cells.Item($r, 18).Value2 = $rVal   # R: Edge total expression weight
    $ws.Cells.Item($r, 19).Value2 = $oVal   # S: Edge average expression derived specificity
    $ws.Cells.Item($r, 20).Value2 = $pVal   # T: Edge total expression derived specificity
}
